$wb = $excel.ActiveWorkbook

$oldGuid = "14d30e09-5902-4338-a1a8-2cf0d70ead89"
$newGuid = "e330d81a-510b-4216-83cc-115e9c78cd9d"

$newFileName    = "$newGuid.md"
$newPathAndName = "e2e\$newGuid.md"
$newHandoffDate = "2016-08-13 15:14:17"

$newZhXlf             = "$newGuid.d397f61b8298cc194c4c3ee166dae59c9472bd8f.zh-cn.xlf"
$newZhHandoffDateTime  = "2016-08-13 15:14:08"
$newHandbackDateTime   = "0001-01-01 00:00:00"

$newDeXlf = "$newGuid.d397f61b8298cc194c4c3ee166dae59c9472bd8f.de-de.xlf"

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("B2").Value = $newPathAndName
$wsOverview.Range("G2").Value = $newHandoffDate
foreach ($h in $wsOverview.Hyperlinks) {
    $h.TextToDisplay = $newPathAndName
}

# --- Sheet "zh-cn" ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newFileName
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDateTime
$wsZh.Range("K2").Value = $newHandbackDateTime

foreach ($h in @($wsZh.Hyperlinks)) {
    if ($h.Range.Address() -eq '$I$2') {
        $h.Delete()
    }
}
foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = $newFileName
    }
}
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsZh.Columns.Item(10).ColumnWidth = 20.833333333333332

# --- Sheet "de-de" ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newFileName
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newHandoffDate
$wsDe.Range("K2").Value = $newHandbackDateTime

foreach ($h in @($wsDe.Hyperlinks)) {
    if ($h.Range.Address() -eq '$I$2') {
        $h.Delete()
    }
}
foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = $newFileName
    }
}
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsDe.Columns.Item(10).ColumnWidth = 20.833333333333332
